$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$data = New-Object 'object[,]' 50,4
$data[0,0] = 0.1151670143008232
$data[0,1] = 0.9680072069168091
$data[0,2] = 0.1076900437474251
$data[0,3] = 0.9787644743919373
$data[1,0] = 0.03377153724431992
$data[1,1] = 0.989281177520752
$data[1,2] = 0.03005184419453144
$data[1,3] = 0.9964147806167603
$data[2,0] = 0.02551163733005524
$data[2,1] = 0.989976704120636
$data[2,2] = 0.02278888784348965
$data[2,3] = 0.9964147806167603
$data[3,0] = 0.02049458771944046
$data[3,1] = 0.9904676079750061
$data[3,2] = 0.0222640298306942
$data[3,3] = 0.9975179433822632
$data[4,0] = 0.01360207796096802
$data[4,1] = 0.9920427203178406
$data[4,2] = 0.0151344295591116
$data[4,3] = 0.9986210465431213
$data[5,0] = 0.01154210232198238
$data[5,1] = 0.9952951669692993
$data[5,2] = 0.01293700188398361
$data[5,3] = 0.9991726279258728
$data[6,0] = 0.01202597934752703
$data[6,1] = 0.994640588760376
$data[6,2] = 0.01113915164023638
$data[6,3] = 0.9991726279258728
$data[7,0] = 0.01160701457411051
$data[7,1] = 0.9944360256195068
$data[7,2] = 0.0103907473385334
$data[7,3] = 0.9991726279258728
$data[8,0] = 0.01071163639426231
$data[8,1] = 0.9949678778648376
$data[8,2] = 0.02177252061665058
$data[8,3] = 0.9983452558517456
$data[9,0] = 0.01132739800959826
$data[9,1] = 0.9949883222579956
$data[9,2] = 0.01582029275596142
$data[9,3] = 0.9983452558517456
$data[10,0] = 0.01048427820205688
$data[10,1] = 0.9951110482215881
$data[10,2] = 0.01228185556828976
$data[10,3] = 0.9991726279258728
$data[11,0] = 0.0112099964171648
$data[11,1] = 0.9944974184036255
$data[11,2] = 0.01321412995457649
$data[11,3] = 0.9991726279258728
$data[12,0] = 0.01039063464850187
$data[12,1] = 0.9949883222579956
$data[12,2] = 0.00936440285295248
$data[12,3] = 0.9991726279258728
$data[13,0] = 0.01039623189717531
$data[13,1] = 0.9949269890785217
$data[13,2] = 0.01437363214790821
$data[13,3] = 0.9983452558517456
$data[14,0] = 0.01045603770762682
$data[14,1] = 0.9949474334716797
$data[14,2] = 0.01627470925450325
$data[14,3] = 0.9983452558517456
$data[15,0] = 0.01105298940092325
$data[15,1] = 0.9949474334716797
$data[15,2] = 0.01197761856019497
$data[15,3] = 0.9991726279258728
$data[16,0] = 0.01121274102479219
$data[16,1] = 0.9945178627967834
$data[16,2] = 0.01125524193048477
$data[16,3] = 0.9991726279258728
$data[17,0] = 0.01078939530998468
$data[17,1] = 0.9949883222579956
$data[17,2] = 0.01310808397829533
$data[17,3] = 0.9991726279258728
$data[18,0] = 0.01093446556478739
$data[18,1] = 0.994640588760376
$data[18,2] = 0.01511598285287619
$data[18,3] = 0.9980695247650146
$data[19,0] = 0.0107086505740881
$data[19,1] = 0.9946610331535339
$data[19,2] = 0.0160983894020319
$data[19,3] = 0.9983452558517456
$data[20,0] = 0.01203931495547295
$data[20,1] = 0.9945791959762573
$data[20,2] = 0.01238434761762619
$data[20,3] = 0.9986210465431213
$data[21,0] = 0.01100021135061979
$data[21,1] = 0.9945996999740601
$data[21,2] = 0.0139040919020772
$data[21,3] = 0.9986210465431213
$data[22,0] = 0.01091726124286652
$data[22,1] = 0.9946610331535339
$data[22,2] = 0.009983655996620655
$data[22,3] = 0.9986210465431213
$data[23,0] = 0.01047796104103327
$data[23,1] = 0.994640588760376
$data[23,2] = 0.01360619533807039
$data[23,3] = 0.9986210465431213
$data[24,0] = 0.01035511400550604
$data[24,1] = 0.9948451519012451
$data[24,2] = 0.007459980901330709
$data[24,3] = 0.9988968372344971
$data[25,0] = 0.0109701044857502
$data[25,1] = 0.9954178929328918
$data[25,2] = 0.01117826718837023
$data[25,3] = 0.9988968372344971
$data[26,0] = 0.009700404480099678
$data[26,1] = 0.995602011680603
$data[26,2] = 0.01277806051075459
$data[26,3] = 0.9991726279258728
$data[27,0] = 0.01026285719126463
$data[27,1] = 0.9950292706489563
$data[27,2] = 0.009460967965424061
$data[27,3] = 0.9986210465431213
$data[28,0] = 0.01055009663105011
$data[28,1] = 0.9947837591171265
$data[28,2] = 0.01394523587077856
$data[28,3] = 0.9988968372344971
$data[29,0] = 0.009984692558646202
$data[29,1] = 0.9952542781829834
$data[29,2] = 0.01364957075566053
$data[29,3] = 0.9986210465431213
$data[30,0] = 0.009636911563575268
$data[30,1] = 0.9950292706489563
$data[30,2] = 0.01593679748475552
$data[30,3] = 0.9991726279258728
$data[31,0] = 0.009974710643291473
$data[31,1] = 0.9953156113624573
$data[31,2] = 0.01925434172153473
$data[31,3] = 0.9980695247650146
$data[32,0] = 0.01019893866032362
$data[32,1] = 0.9948451519012451
$data[32,2] = 0.01827728562057018
$data[32,3] = 0.9991726279258728
$data[33,0] = 0.01052526105195284
$data[33,1] = 0.9951315522193909
$data[33,2] = 0.01648190803825855
$data[33,3] = 0.9986210465431213
$data[34,0] = 0.0113256424665451
$data[34,1] = 0.9944769740104675
$data[34,2] = 0.01323280856013298
$data[34,3] = 0.9991726279258728
$data[35,0] = 0.009807305410504341
$data[35,1] = 0.9953156113624573
$data[35,2] = 0.01340371277183294
$data[35,3] = 0.9991726279258728
$data[36,0] = 0.01031882874667645
$data[36,1] = 0.9947837591171265
$data[36,2] = 0.01601662673056126
$data[36,3] = 0.9988968372344971
$data[37,0] = 0.009958542883396149
$data[37,1] = 0.9951315522193909
$data[37,2] = 0.01539130602031946
$data[37,3] = 0.9988968372344971
$data[38,0] = 0.009725025855004787
$data[38,1] = 0.995602011680603
$data[38,2] = 0.01417617592960596
$data[38,3] = 0.9991726279258728
$data[39,0] = 0.010613146238029
$data[39,1] = 0.994886040687561
$data[39,2] = 0.01684320531785488
$data[39,3] = 0.9991726279258728
$data[40,0] = 0.01019390113651752
$data[40,1] = 0.994640588760376
$data[40,2] = 0.02538632787764072
$data[40,3] = 0.9986210465431213
$data[41,0] = 0.01104004960507154
$data[41,1] = 0.9949269890785217
$data[41,2] = 0.0166435856372118
$data[41,3] = 0.9991726279258728
$data[42,0] = 0.01068101357668638
$data[42,1] = 0.995356559753418
$data[42,2] = 0.005823639687150717
$data[42,3] = 0.9997242093086243
$data[43,0] = 0.01098092366009951
$data[43,1] = 0.9951315522193909
$data[43,2] = 0.01703402958810329
$data[43,3] = 0.9991726279258728
$data[44,0] = 0.01040267013013363
$data[44,1] = 0.9947019815444946
$data[44,2] = 0.01886876113712788
$data[44,3] = 0.9991726279258728
$data[45,0] = 0.01065365131944418
$data[45,1] = 0.9945996999740601
$data[45,2] = 0.02223456464707851
$data[45,3] = 0.9988968372344971
$data[46,0] = 0.01165830809623003
$data[46,1] = 0.9943132996559143
$data[46,2] = 0.02060576342046261
$data[46,3] = 0.9991726279258728
$data[47,0] = 0.01000626850873232
$data[47,1] = 0.9951928853988647
$data[47,2] = 0.01971776969730854
$data[47,3] = 0.9991726279258728
$data[48,0] = 0.01030280441045761
$data[48,1] = 0.9951315522193909
$data[48,2] = 0.02359041757881641
$data[48,3] = 0.9988968372344971
$data[49,0] = 0.01173094287514687
$data[49,1] = 0.9945996999740601
$data[49,2] = 0.01725848950445652
$data[49,3] = 0.9983452558517456
$ws.Range("A2:D51").Value = $data
Write-Output "done"
